$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("J_Montgomery78@yahoo.com", 26940137),
    @("Caledoni_C@xs4all.nl", 79667400),
    @("Aile.B52@gmail.com", 80613716),
    @("Ma.Chavez70@yahoo.com", 62538332),
    @("JaG@hotmail.com", 42152129),
    @("Alexander-JameHansen@xs4all.nl", 97408675)
)

$startRow = 117
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
